# Fix NA values in the input for Dynare for the ifn job.
# Populates previously-empty numeric cells in rows 2 and 102 (columns B:AT)
# of the active worksheet with the correct values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, columns B2:AT2
$row2 = @(0.56, 0.78, 0.23431, 0.067, 3.1, 0.1, 0.09, 0.18, 0.08, 0.06, 0, $Null, $Null, 0.03, 0.07446, 90.95893, 207, 1, 0.19, 0.7, 0.18, 0, 0.034, 1, 1.42, 1, 1, $Null, 0.065, 0.091, 1.85, 0.091, $Null, $Null, 0.5, $Null, 1, 1, 23, 0.1, 185, 99.0099009901, 267, 185, 0.25)

# Row 102, columns B102:AT102
$row102 = @(0.56, 0.78, 0.1, 0.067, 3.1, 0.1, 0.09, 0.18, 0.08, 0.06, 0, 0, 0.01, 0.03, 0.08, 216, 207, 0.599, 0.19, 0.7, 0.13878, 0, 0.034, 1, 1.42, 1, 1, $Null, 0.055, 0.091, 1.85, 0.091, $Null, 0.4, 0.5, $Null, 1, 1, 23, 0.1, 185, 265, 267, 185, 0.25)

$startCol = 2  # column B
for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = $startCol + $i
    $val = $row2[$i]
    if ($val -ne $Null) {
        $ws.Cells.Item(2, $col).Value2 = $val
    }
}

for ($i = 0; $i -lt $row102.Length; $i++) {
    $col = $startCol + $i
    $val = $row102[$i]
    if ($val -ne $Null) {
        $ws.Cells.Item(102, $col).Value2 = $val
    }
}
